# The presentation currently uses the "Integral" theme (colour scheme
# stored in ppt/theme/theme2.xml, which is what the slide master / all
# slides reference). The commit swaps the content of the two theme parts
# (ppt/theme/theme1.xml <-> ppt/theme/theme2.xml) so that the design that
# is actually applied to the slides becomes the stock "Office Theme"
# colour palette instead of "Integral".
#
# Drive this through the Slide ThemeColorScheme (the DrawingML-backed
# 12-slot clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) so
# every slide (they all share the one slide master / theme part) picks
# up the standard Office theme colours.

$p = $ppt.ActivePresentation

# Office Theme colour scheme (RRGGBB), in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($s = 1; $s -le $p.Slides.Count; $s++) {
    $slide = $p.Slides.Item($s)
    $tcs = $slide.ThemeColorScheme
    for ($i = 1; $i -le 12; $i++) {
        $tcs.Item($i).RGB = $officeThemeRgb[$i - 1]
    }
}
